# Left-align the certificate body paragraphs.
#
# The student-details / dates / host-company / attestation paragraphs
# (which describe the completed internship) are switched from the
# default (justify/none) alignment to explicit left alignment, matching
# the paragraph right above them ("ολοκλήρωσε την Πρακτική Άσκηση...")
# that was already left-aligned.

$d = $word.ActiveDocument

$targets = @(
    "Βεβαιώνεται ότι ο/η",
    "στο χρονικό διάστημα",
    "στον Φορ",
    "Μετά από επικοινωνία με τον φορέα υποδοχής"
)

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    foreach ($needle in $targets) {
        if ($t.StartsWith($needle)) {
            $p.Range.ParagraphFormat.Alignment = 0
            break
        }
    }
}
